$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp column (O) for rows 2 through 72 from
# "2022-08-26 07:07:07" to "2022-08-26 21:00:03"
for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    if ($cell.Value2 -eq "2022-08-26 07:07:07") {
        $cell.Value = "2022-08-26 21:00:03"
    }
}

# Update the productAriaLabel text for row 65 (column N) to reflect
# that the item is out of stock online.
$ws.Range("N65").Value = "Naturaline Herren T-Shirt Kurzarm weiss XXL - Online kein Bestand 24.95 Schweizer Franken"
